$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.150997519493103
$ws.Range("B1").Value = 2.39935564994812
$ws.Range("C1").Value = 5.141182899475098
$ws.Range("D1").Value = 2.212102174758911
$ws.Range("E1").Value = 1.245173573493958
